$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 11 ("2021年") below the existing data (previous last row was 10 / "2020年").
# Copy the row-label cell formatting (bold, bordered, centered) from A10 to A11,
# matching the style already used for every other year label in column A.
$srcLabel = $ws.Cells.Item(10, 1)
$dstLabel = $ws.Cells.Item(11, 1)
$srcLabel.Copy()
$dstLabel.PasteSpecial(-4122)   # xlPasteFormats
$dstLabel.Value = "2021年"

# Numeric indicator values for 2021, one per industry column (B..AQ).
# Column E ("其他采矿业私营工业企业单位数") has no reported value for 2021 and stays blank.
$ws.Range("B11").Value = 17635
$ws.Range("C11").Value = 4155
$ws.Range("D11").Value = 1469
$ws.Range("F11").Value = 17855
$ws.Range("G11").Value = 15947
$ws.Range("H11").Value = 1673
$ws.Range("I11").Value = 5345
$ws.Range("J11").Value = 5150
$ws.Range("K11").Value = 6020
$ws.Range("L11").Value = 1988
$ws.Range("M11").Value = 228
$ws.Range("N11").Value = 7856
$ws.Range("O11").Value = 6414
$ws.Range("P11").Value = 650
$ws.Range("Q11").Value = 10438
$ws.Range("R11").Value = 18233
$ws.Range("S11").Value = 522
$ws.Range("T11").Value = 11211
$ws.Range("U11").Value = 28
$ws.Range("V11").Value = 2805
$ws.Range("W11").Value = 934
$ws.Range("X11").Value = 2375
$ws.Range("Y11").Value = 22901
$ws.Range("Z11").Value = 6978
$ws.Range("AA11").Value = 1509
$ws.Range("AB11").Value = 16
$ws.Range("AC11").Value = 325752
$ws.Range("AD11").Value = 16964
$ws.Range("AE11").Value = 10311
$ws.Range("AF11").Value = 15094
$ws.Range("AG11").Value = 23505
$ws.Range("AH11").Value = 5727
$ws.Range("AI11").Value = 3950
$ws.Range("AJ11").Value = 313
$ws.Range("AK11").Value = 25430
$ws.Range("AL11").Value = 3924
$ws.Range("AM11").Value = 35577
$ws.Range("AN11").Value = 2743
$ws.Range("AO11").Value = 6297
$ws.Range("AP11").Value = 4472
$ws.Range("AQ11").Value = 1108
